$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$transitionProbabilities = @{
    "B2" = 0.2666666666666667
    "C2" = 0.4
    "P2" = 0.1333333333333333
    "S2" = 0.2
    "P3" = 0.5
    "S3" = 0.5
    "S4" = 1
    "B6" = 0.1333333333333333
    "F6" = 0.06666666666666667
    "J6" = 0.1333333333333333
    "Q6" = 0.2
    "R6" = 0.06666666666666667
    "S6" = 0.4
    "F7" = 0.1
    "O7" = 0.1
    "Q7" = 0.3
    "S7" = 0.5
    "F8" = 0.1282051282051282
    "J8" = 0.2051282051282051
    "Q8" = 0.2307692307692308
    "R8" = 0.1025641025641026
    "S8" = 0.3333333333333333
    "F9" = 0.06666666666666667
    "J9" = 0.2
    "Q9" = 0.06666666666666667
    "R9" = 0.06666666666666667
    "S9" = 0.6
    "B10" = 0.09677419354838709
    "D10" = 0.02150537634408602
    "F10" = 0.02150537634408602
    "J10" = 0.09677419354838709
    "O10" = 0.03225806451612903
    "Q10" = 0.3010752688172043
    "R10" = 0.03225806451612903
    "S10" = 0.3978494623655914
    "J11" = 0.1875
    "K11" = 0.0625
    "L11" = 0.75
    "G12" = 0.75
    "J12" = 0.25
    "G13" = 1
    "H15" = 0.1875
    "I15" = 0.125
    "J15" = 0.375
    "K15" = 0.0625
    "S15" = 0.25
    "F16" = 0.2
    "H16" = 0.2
    "J16" = 0.4
    "S16" = 0.2
    "H17" = 0.2954545454545455
    "I17" = 0.04545454545454546
    "J17" = 0.4772727272727273
    "K17" = 0.06818181818181818
    "O17" = 0.06818181818181818
    "S17" = 0.04545454545454546
    "H18" = 0.2222222222222222
    "I18" = 0.1111111111111111
    "J18" = 0.3333333333333333
    "K18" = 0.2222222222222222
    "O18" = 0.1111111111111111
    "F19" = 0.02272727272727273
    "H19" = 0.2272727272727273
    "I19" = 0.1136363636363636
    "J19" = 0.3863636363636364
    "K19" = 0.1022727272727273
    "M19" = 0.02272727272727273
    "O19" = 0.07954545454545454
    "S19" = 0.04545454545454546
}

foreach ($cellRef in $transitionProbabilities.Keys) {
    $ws.Range($cellRef).Value = $transitionProbabilities[$cellRef]
}

